$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.574.24"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.840.86"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -2.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4301"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3728"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07287"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8704"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "1.875.14"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.389"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "27.602.09"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.185"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").Value = "2.069.23"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  -4.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +7.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.313"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08899"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.214"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7721"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.511"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.902"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.008"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.125"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01968"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05292"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.883"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.120"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1685"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.722"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4735"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06438"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.008"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.681"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("E51").Value = "  -2.80%  "
